# Update the LABELS column (E) values: each numeric label cycles forward by one,
# wrapping modulo 4 (0->1, 1->2, 2->3, 3->0), for data rows 2 through 476.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 476
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    $current = $cell.Value2
    if ($current -ne $null) {
        $new = ($current + 1) % 4
        $cell.Value = $new
    }
}
